$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K" header) values for rows 2-13
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 5
$ws.Range("G4").Value = 6
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 4
$ws.Range("G7").Value = 3
$ws.Range("G8").Value = 3
$ws.Range("G9").Value = 5
$ws.Range("G10").Value = 3
$ws.Range("G11").Value = 1
$ws.Range("G12").Value = 0
$ws.Range("G13").Value = 1
